$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the new DP problem entry under "DP on Matrix" column (E)
$ws.Range("E5").Value = "Leetcode - 542"

# Update the active selection to match the new edit location
$ws.Range("E8").Select()
